$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 27.07648166666667
$ws.Range("H2").Value = 81.229445
$ws.Range("I2").Value = 0.1943552322922666
$ws.Range("J2").Value = 0.1943552322922666
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.580817
$ws.Range("N2").Value = 1.742451
$ws.Range("O2").Value = 0.1705489461834183
$ws.Range("P2").Value = 0.1705489461834183
$ws.Range("Q2").Value = 15.72648085218833
$ws.Range("R2").Value = 141.538327669695
$ws.Range("S2").Value = 0.03314708005267955
$ws.Range("T2").Value = 0.03314708005267955
$ws.Range("G3").Value = 27.07648166666667
$ws.Range("H3").Value = 81.229445
$ws.Range("I3").Value = 0.1943552322922666
$ws.Range("J3").Value = 0.1943552322922666
$ws.Range("O3").Value = 0.6679715536912479
$ws.Range("P3").Value = 0.6679715536912479
$ws.Range("Q3").Value = 61.59429350935056
$ws.Range("R3").Value = 554.348641584155
$ws.Range("S3").Value = 0.1298237664822887
$ws.Range("T3").Value = 0.1298237664822887
$ws.Range("G4").Value = 27.07648166666667
$ws.Range("H4").Value = 81.229445
$ws.Range("I4").Value = 0.1943552322922666
$ws.Range("J4").Value = 0.1943552322922666
$ws.Range("M4").Value = 0.5499303333333333
$ws.Range("N4").Value = 1.649791
$ws.Range("O4").Value = 0.1614795001253337
$ws.Range("P4").Value = 0.1614795001253338
$ws.Range("Q4").Value = 14.89017858844389
$ws.Range("R4").Value = 134.011607295995
$ws.Range("S4").Value = 0.03138438575729834
$ws.Range("T4").Value = 0.03138438575729834
$ws.Range("I5").Value = 0.599012687336886
$ws.Range("J5").Value = 0.599012687336886
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.580817
$ws.Range("N5").Value = 1.742451
$ws.Range("O5").Value = 0.1705489461834183
$ws.Range("P5").Value = 0.1705489461834183
$ws.Range("Q5").Value = 48.46981193413567
$ws.Range("R5").Value = 436.228307407221
$ws.Range("S5").Value = 0.1021609825758034
$ws.Range("T5").Value = 0.1021609825758034
$ws.Range("I6").Value = 0.599012687336886
$ws.Range("J6").Value = 0.599012687336886
$ws.Range("O6").Value = 0.6679715536912479
$ws.Range("P6").Value = 0.6679715536912479
$ws.Range("S6").Value = 0.4001234354411894
$ws.Range("T6").Value = 0.4001234354411894
$ws.Range("I7").Value = 0.599012687336886
$ws.Range("J7").Value = 0.599012687336886
$ws.Range("M7").Value = 0.5499303333333333
$ws.Range("N7").Value = 1.649791
$ws.Range("O7").Value = 0.1614795001253337
$ws.Range("P7").Value = 0.1614795001253338
$ws.Range("Q7").Value = 45.89228592404011
$ws.Range("R7").Value = 413.030573316361
$ws.Range("S7").Value = 0.09672826931989317
$ws.Range("T7").Value = 0.09672826931989319
$ws.Range("G8").Value = 14.445417
$ws.Range("H8").Value = 43.336251
$ws.Range("I8").Value = 0.1036893349422856
$ws.Range("J8").Value = 0.1036893349422857
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.580817
$ws.Range("N8").Value = 1.742451
$ws.Range("O8").Value = 0.1705489461834183
$ws.Range("P8").Value = 0.1705489461834183
$ws.Range("Q8").Value = 8.390143765689
$ws.Range("R8").Value = 75.51129389120099
$ws.Range("S8").Value = 0.01768410680486631
$ws.Range("T8").Value = 0.01768410680486632
$ws.Range("G9").Value = 14.445417
$ws.Range("H9").Value = 43.336251
$ws.Range("I9").Value = 0.1036893349422856
$ws.Range("J9").Value = 0.1036893349422857
$ws.Range("O9").Value = 0.6679715536912479
$ws.Range("P9").Value = 0.6679715536912479
$ws.Range("Q9").Value = 32.860814987581
$ws.Range("R9").Value = 295.747334888229
$ws.Range("S9").Value = 0.06926152616261075
$ws.Range("T9").Value = 0.06926152616261076
$ws.Range("G10").Value = 14.445417
$ws.Range("H10").Value = 43.336251
$ws.Range("I10").Value = 0.1036893349422856
$ws.Range("J10").Value = 0.1036893349422857
$ws.Range("M10").Value = 0.5499303333333333
$ws.Range("N10").Value = 1.649791
$ws.Range("O10").Value = 0.1614795001253337
$ws.Range("P10").Value = 0.1614795001253338
$ws.Range("Q10").Value = 7.943972985948999
$ws.Range("R10").Value = 71.495756873541
$ws.Range("S10").Value = 0.01674370197480859
$ws.Range("T10").Value = 0.01674370197480859
$ws.Range("G11").Value = 14.34140633333333
$ws.Range("H11").Value = 43.024219
$ws.Range("I11").Value = 0.1029427454285617
$ws.Range("J11").Value = 0.1029427454285617
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.580817
$ws.Range("N11").Value = 1.742451
$ws.Range("O11").Value = 0.1705489461834183
$ws.Range("P11").Value = 0.1705489461834183
$ws.Range("Q11").Value = 8.329732602307667
$ws.Range("R11").Value = 74.96759342076901
$ws.Range("S11").Value = 0.0175567767500691
$ws.Range("T11").Value = 0.0175567767500691
$ws.Range("G12").Value = 14.34140633333333
$ws.Range("H12").Value = 43.024219
$ws.Range("I12").Value = 0.1029427454285617
$ws.Range("J12").Value = 0.1029427454285617
$ws.Range("O12").Value = 0.6679715536912479
$ws.Range("P12").Value = 0.6679715536912479
$ws.Range("Q12").Value = 32.62420878410011
$ws.Range("R12").Value = 293.617879056901
$ws.Range("S12").Value = 0.06876282560515895
$ws.Range("T12").Value = 0.06876282560515895
$ws.Range("G13").Value = 14.34140633333333
$ws.Range("H13").Value = 43.024219
$ws.Range("I13").Value = 0.1029427454285617
$ws.Range("J13").Value = 0.1029427454285617
$ws.Range("M13").Value = 0.5499303333333333
$ws.Range("N13").Value = 1.649791
$ws.Range("O13").Value = 0.1614795001253337
$ws.Range("P13").Value = 0.1614795001253338
$ws.Range("Q13").Value = 7.886774365358777
$ws.Range("R13").Value = 70.98096928822901
$ws.Range("S13").Value = 0.01662314307333362
$ws.Range("T13").Value = 0.01662314307333363
